$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right
#    after the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaPara  = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text.Contains("Meta description")) {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2) At the very end of the document, insert a new bold paragraph
#    ("Play 777 Heist Free: ...") right before the final paragraph,
#    and replace the final paragraph's italic text with the new
#    "Read our review ..." copy (keeping the italic formatting).
# ------------------------------------------------------------------
$count     = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)

$boldText   = "Play 777 Heist Free: Impressive Graphics and Bonus Features"
$italicText = "Read our review of 777 Heist, an online slot game with impressive graphics, bonus features, and numerous opportunities for big payouts. Play 777 Heist free."

$targetRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p>' +
            '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $italicText + '</w:t></w:r></w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

[void]$targetRange.InsertXML($xmlFragment)
